$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("98, 130, 786, 786", 1, "786", "594"),
    @("130, 426, 458", 1, "458", "692"),
    @("130, 130, 458", 1, "130", "565"),
    @("130, 130, 426, 786", 1, "130", "969"),
    @("130, 458, 754, 786", 1, "786", "707"),
    @("130, 786, 1082", 1, "130", "787"),
    @("98, 130, 786, SF", 1, "786", "959"),
    @("130, 1082, 1114, SF", 1, "1114", "777"),
    @("426, 754, 1082, SF", 4, "754, SF, 1082, 426", "1248, 1248, 1248, 1248"),
    @("98, 130, 786, 1082", 1, "786", "794"),
    @("98, 426, 1082, SF, SF", 1, "SF", "806"),
    @("98, 458, 754, 786, 1082", 1, "1082", "946"),
    @("98, 458, 754, 1082", 1, "458", "1186"),
    @("130, 426, 786, SF, SF", 1, "786", "1225"),
    @("98, 130, 754, 786, 1082, SF", 1, "130", "1219")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $data[$i][0]

    $ws.Cells.Item($row, 2).Value = $data[$i][1]

    $cC = $ws.Cells.Item($row, 3)
    $cC.NumberFormat = "@"
    $cC.Value = $data[$i][2]
    $cC.ClearFormats()

    $cD = $ws.Cells.Item($row, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $data[$i][3]
    $cD.ClearFormats()
}
